$wb = $excel.ActiveWorkbook

# --- DCpUC sheet: repurpose "geothermal" row's cost to reference "hydro" (B5) instead of the
#     onshore-wind based Data!C6 lookup, and drop the now-unneeded yellow highlight fill.
$wsDCpUC = $wb.Worksheets.Item("DCpUC")
$wsDCpUC.Range("B10").Interior.Pattern = -4142
$wsDCpUC.Range("B10").Formula = "=B5"
$wsDCpUC.Range("B10").NumberFormat = "0"

# --- About sheet: add an explanatory note about the geothermal -> pumped hydro repurposing.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A30").Value = "In the India EPS, the geothermal plant type is repurposed as pumped hydro capacity."
$wsAbout.Range("A30").Font.Color = 0
$wsAbout.Range("A30").VerticalAlignment = -4108
